$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "Player Info" worksheet as the first sheet -----------
# Worksheets.Add() inserts a new sheet before the currently active sheet,
# which is exactly where the diff places "Player Info" (ahead of the two
# existing "ODI ..." sheets).
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $playerInfo.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the bold / bordered / centered-top header styling used by the
# other two sheets' header rows.
$headerRow = $playerInfo.Range("A1:D1")
$headerRow.Font.Bold = $true
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160
$headerRow.Borders.LineStyle = 1

# Player data row. The leading "'" keeps the numeric-looking ID as text,
# matching the source data (every cell in this workbook is stored as text).
$playerInfo.Range("A2").Value = "'4312"
$playerInfo.Range("B2").Value = "Jermaine Blackwood"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# --- 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").Value = "'3853"
$batting.Range("D3").Value = "'3855"
$batting.Range("D4").Value = "'4636"

# --- 3. "ODI Bowling": same rename + conversion ----------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = "'3855"
